$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying shared-string table order changed (re-run of the notebook
# produced a different, but equivalent, ordering for goods with tied counts).
# Column B (counts) stays aligned to its row; only the word in column A changes
# for the rows whose shared string was reshuffled.
$ws.Range("A16").Value = "колеса"
$ws.Range("A18").Value = "полотно"
$ws.Range("A19").Value = "говядина"
$ws.Range("A21").Value = "парча"
$ws.Range("A22").Value = "позумент"
$ws.Range("A23").Value = "табак"
$ws.Range("A26").Value = "сахар"
$ws.Range("A27").Value = "выбойка"
$ws.Range("A30").Value = "ладан"
$ws.Range("A32").Value = "китайка"
$ws.Range("A33").Value = "сапог"
$ws.Range("A35").Value = "конь"
$ws.Range("A36").Value = "платок"
$ws.Range("A37").Value = "рогожа"
$ws.Range("A38").Value = "замок"
$ws.Range("A39").Value = "гвоздь"
$ws.Range("A40").Value = "овца"
$ws.Range("A41").Value = "обод"
$ws.Range("A42").Value = "горшок"
$ws.Range("A43").Value = "веревка"
$ws.Range("A44").Value = "ром"
$ws.Range("A45").Value = "котел"
$ws.Range("A46").Value = "гумми"
$ws.Range("A47").Value = "хомут"
$ws.Range("A48").Value = "брусья"
$ws.Range("A49").Value = "нитка"
$ws.Range("A50").Value = "роза"
$ws.Range("A51").Value = "дуга"
$ws.Range("A52").Value = "сосуд"
$ws.Range("A53").Value = "скотский кожа"
$ws.Range("A54").Value = "бечева"
$ws.Range("A55").Value = "покроми"
$ws.Range("A56").Value = "сковорода"
